{"js": "// Apply the set of before->after text replacements recorded in the diff.\n// Every \"before\" string occurs exactly once in the document (verified against\n// the source OOXML), so an exact-text search + single-result replace is safe and\n// unambiguous; pairs are applied in document order, matching the diff.\nconst pairs = [\n  [\"2023-07-02 Sunday\", \"2023-07-03 Monday\"],\n  [\"19-16=3\", \"18+9=27\"],\n  [\"34+15=49\", \"46-33=13\"],\n  [\"34+12=46\", \"87-87=0\"],\n  [\"60+39=99\", \"18-11=7\"],\n  [\"12+69=81\", \"49+16=65\"],\n  [\"93-25=68\", \"43+54=97\"],\n  [\"18+23=41\", \"76-52=24\"],\n  [\"1+93=94\", \"91-48=43\"],\n  [\"2+50=52\", \"35-25=10\"],\n  [\"40-22=18\", \"69-36=33\"],\n  [\"2+68=70\", \"27+34=61\"],\n  [\"45-24=21\", \"46+53=99\"],\n  [\"52+1=53\", \"22+15=37\"],\n  [\"21+5=26\", \"6+41=47\"],\n  [\"73-62=11\", \"58-48=10\"],\n  [\"73-66=7\", \"45-10=35\"],\n  [\"90-8=82\", \"36+2=38\"],\n  [\"9+46=55\", \"16+81=97\"],\n  [\"30-15=15\", \"20+16=36\"],\n  [\"63+11=74\", \"76-45=31\"],\n  [\"82-3=79\", \"64-14=50\"],\n  [\"28+58=86\", \"8+38=46\"],\n  [\"96-71=25\", \"11+22=33\"],\n  [\"23+44=67\", \"77-50=27\"],\n  [\"4+19=23\", \"2+24=26\"],\n  [\"9+70=79\", \"50+48=98\"],\n  [\"57-16=41\", \"96-80=16\"],\n  [\"64-9=55\", \"95-88=7\"],\n  [\"0+4=4\", \"9+39=48\"],\n  [\"35+47=82\", \"22-2=20\"],\n  [\"0+8=8\", \"68-24=44\"],\n  [\"41+3=44\", \"40+4=44\"],\n  [\"78-4=74\", \"95-46=49\"],\n  [\"39-13=26\", \"32+43=75\"],\n  [\"62+20=82\", \"19+8=27\"],\n  [\"2+92=94\", \"5+85=90\"],\n  [\"32+23=55\", \"73-0=73\"],\n  [\"86+0=86\", \"68+11=79\"],\n  [\"4+61=65\", \"95-74=21\"],\n  [\"34+23=57\", \"14+13=27\"],\n  [\"39-2=37\", \"56-19=37\"],\n  [\"8+54=62\", \"89+9=98\"],\n  [\"44+27=71\", \"34+51=85\"],\n  [\"26-21=5\", \"95-11=84\"],\n  [\"96-72=24\", \"25+27=52\"],\n  [\"53+39=92\", \"83-47=36\"],\n  [\"21+61=82\", \"27-20=7\"],\n  [\"42-24=18\", \"0+54=54\"],\n  [\"62-37=25\", \"51-49=2\"],\n  [\"38-38=0\", \"81+3=84\"],\n  [\"62+15=77\", \"90-5=85\"],\n  [\"1+88=89\", \"20+42=62\"],\n  [\"12-4=8\", \"87-18=69\"],\n  [\"97-63=34\", \"69-44=25\"],\n  [\"89-24=65\", \"76-64=12\"],\n  [\"81-54=27\", \"85+5=90\"],\n  [\"40-7=33\", \"25-24=1\"],\n  [\"22+60=82\", \"62-20=42\"],\n  [\"62-51=11\", \"96-50=46\"],\n  [\"63-22=41\", \"97-28=69\"],\n  [\"67-1=66\", \"27-26=1\"],\n  [\"19+3=22\", \"10+69=79\"],\n  [\"26+37=63\", \"31+67=98\"],\n  [\"21-5=16\", \"1+13=14\"],\n  [\"96-69=27\", \"79-19=60\"],\n  [\"10+8=18\", \"80-79=1\"],\n  [\"36+63=99\", \"2+13=15\"],\n  [\"53-10=43\", \"3+54=57\"],\n  [\"97-4=93\", \"9+10=19\"],\n  [\"37-16=21\", \"2+68=70\"],\n  [\"74-9=65\", \"84-37=47\"],\n  [\"72+26=98\", \"59-21=38\"],\n  [\"18+24=42\", \"52-9=43\"],\n  [\"16+67=83\", \"10+88=98\"],\n  [\"55+17=72\", \"54-14=40\"],\n  [\"50-27=23\", \"48-40=8\"],\n  [\"39-25=14\", \"78+11=89\"],\n  [\"86-73=13\", \"1+74=75\"],\n  [\"10+64=74\", \"19-15=4\"],\n  [\"85-33=52\", \"97-47=50\"],\n  [\"44-26=18\", \"41+54=95\"],\n  [\"5+72=77\", \"30-23=7\"],\n  [\"25-9=16\", \"97-68=29\"],\n  [\"50-20=30\", \"5+43=48\"],\n  [\"97-1=96\", \"74-49=25\"],\n  [\"84-56=28\", \"40+18=58\"],\n  [\"74-51=23\", \"78-49=29\"],\n  [\"45-13=32\", \"67+23=90\"],\n  [\"2+2=4\", \"5+54=59\"],\n  [\"89-64=25\", \"50+28=78\"],\n  [\"48+51=99\", \"35+37=72\"],\n  [\"56-10=46\", \"8+10=18\"],\n  [\"71+6=77\", \"44-9=35\"],\n  [\"33+14=47\", \"65+28=93\"],\n  [\"24+59=83\", \"47-33=14\"],\n  [\"30+69=99\", \"98-83=15\"],\n  [\"35+43=78\", \"0+56=56\"],\n  [\"47-32=15\", \"1+6=7\"],\n  [\"65-8=57\", \"14+77=91\"],\n  [\"25+60=85\", \"49+14=63\"],\n];\n\nfor (const [before, after] of pairs) {\n  const results = context.document.body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      `Expected exactly one match for ${JSON.stringify(before)}, found ${results.items.length}`\n    );\n  }\n\n  results.items[0].insertText(after, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Apply the set of before->after text replacements recorded in the diff.\n# Each \"before\" string is unique in the document, so a literal Find/Replace\n# (wdReplaceOne) on the whole document Content range is unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2023-07-02 Sunday\", \"2023-07-03 Monday\"),\n    @(\"19-16=3\", \"18+9=27\"),\n    @(\"34+15=49\", \"46-33=13\"),\n    @(\"34+12=46\", \"87-87=0\"),\n    @(\"60+39=99\", \"18-11=7\"),\n    @(\"12+69=81\", \"49+16=65\"),\n    @(\"93-25=68\", \"43+54=97\"),\n    @(\"18+23=41\", \"76-52=24\"),\n    @(\"1+93=94\", \"91-48=43\"),\n    @(\"2+50=52\", \"35-25=10\"),\n    @(\"40-22=18\", \"69-36=33\"),\n    @(\"2+68=70\", \"27+34=61\"),\n    @(\"45-24=21\", \"46+53=99\"),\n    @(\"52+1=53\", \"22+15=37\"),\n    @(\"21+5=26\", \"6+41=47\"),\n    @(\"73-62=11\", \"58-48=10\"),\n    @(\"73-66=7\", \"45-10=35\"),\n    @(\"90-8=82\", \"36+2=38\"),\n    @(\"9+46=55\", \"16+81=97\"),\n    @(\"30-15=15\", \"20+16=36\"),\n    @(\"63+11=74\", \"76-45=31\"),\n    @(\"82-3=79\", \"64-14=50\"),\n    @(\"28+58=86\", \"8+38=46\"),\n    @(\"96-71=25\", \"11+22=33\"),\n    @(\"23+44=67\", \"77-50=27\"),\n    @(\"4+19=23\", \"2+24=26\"),\n    @(\"9+70=79\", \"50+48=98\"),\n    @(\"57-16=41\", \"96-80=16\"),\n    @(\"64-9=55\", \"95-88=7\"),\n    @(\"0+4=4\", \"9+39=48\"),\n    @(\"35+47=82\", \"22-2=20\"),\n    @(\"0+8=8\", \"68-24=44\"),\n    @(\"41+3=44\", \"40+4=44\"),\n    @(\"78-4=74\", \"95-46=49\"),\n    @(\"39-13=26\", \"32+43=75\"),\n    @(\"62+20=82\", \"19+8=27\"),\n    @(\"2+92=94\", \"5+85=90\"),\n    @(\"32+23=55\", \"73-0=73\"),\n    @(\"86+0=86\", \"68+11=79\"),\n    @(\"4+61=65\", \"95-74=21\"),\n    @(\"34+23=57\", \"14+13=27\"),\n    @(\"39-2=37\", \"56-19=37\"),\n    @(\"8+54=62\", \"89+9=98\"),\n    @(\"44+27=71\", \"34+51=85\"),\n    @(\"26-21=5\", \"95-11=84\"),\n    @(\"96-72=24\", \"25+27=52\"),\n    @(\"53+39=92\", \"83-47=36\"),\n    @(\"21+61=82\", \"27-20=7\"),\n    @(\"42-24=18\", \"0+54=54\"),\n    @(\"62-37=25\", \"51-49=2\"),\n    @(\"38-38=0\", \"81+3=84\"),\n    @(\"62+15=77\", \"90-5=85\"),\n    @(\"1+88=89\", \"20+42=62\"),\n    @(\"12-4=8\", \"87-18=69\"),\n    @(\"97-63=34\", \"69-44=25\"),\n    @(\"89-24=65\", \"76-64=12\"),\n    @(\"81-54=27\", \"85+5=90\"),\n    @(\"40-7=33\", \"25-24=1\"),\n    @(\"22+60=82\", \"62-20=42\"),\n    @(\"62-51=11\", \"96-50=46\"),\n    @(\"63-22=41\", \"97-28=69\"),\n    @(\"67-1=66\", \"27-26=1\"),\n    @(\"19+3=22\", \"10+69=79\"),\n    @(\"26+37=63\", \"31+67=98\"),\n    @(\"21-5=16\", \"1+13=14\"),\n    @(\"96-69=27\", \"79-19=60\"),\n    @(\"10+8=18\", \"80-79=1\"),\n    @(\"36+63=99\", \"2+13=15\"),\n    @(\"53-10=43\", \"3+54=57\"),\n    @(\"97-4=93\", \"9+10=19\"),\n    @(\"37-16=21\", \"2+68=70\"),\n    @(\"74-9=65\", \"84-37=47\"),\n    @(\"72+26=98\", \"59-21=38\"),\n    @(\"18+24=42\", \"52-9=43\"),\n    @(\"16+67=83\", \"10+88=98\"),\n    @(\"55+17=72\", \"54-14=40\"),\n    @(\"50-27=23\", \"48-40=8\"),\n    @(\"39-25=14\", \"78+11=89\"),\n    @(\"86-73=13\", \"1+74=75\"),\n    @(\"10+64=74\", \"19-15=4\"),\n    @(\"85-33=52\", \"97-47=50\"),\n    @(\"44-26=18\", \"41+54=95\"),\n    @(\"5+72=77\", \"30-23=7\"),\n    @(\"25-9=16\", \"97-68=29\"),\n    @(\"50-20=30\", \"5+43=48\"),\n    @(\"97-1=96\", \"74-49=25\"),\n    @(\"84-56=28\", \"40+18=58\"),\n    @(\"74-51=23\", \"78-49=29\"),\n    @(\"45-13=32\", \"67+23=90\"),\n    @(\"2+2=4\", \"5+54=59\"),\n    @(\"89-64=25\", \"50+28=78\"),\n    @(\"48+51=99\", \"35+37=72\"),\n    @(\"56-10=46\", \"8+10=18\"),\n    @(\"71+6=77\", \"44-9=35\"),\n    @(\"33+14=47\", \"65+28=93\"),\n    @(\"24+59=83\", \"47-33=14\"),\n    @(\"30+69=99\", \"98-83=15\"),\n    @(\"35+43=78\", \"0+56=56\"),\n    @(\"47-32=15\", \"1+6=7\"),\n    @(\"65-8=57\", \"14+77=91\"),\n    @(\"25+60=85\", \"49+14=63\")\n)\n\nforeach ($pair in $pairs) {\n    $before = $pair[0]\n    $after = $pair[1]\n\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n\n    # Wrap = wdFindContinue (1); Replace = wdReplaceOne (1) -- only one\n    # occurrence is expected/required for each \"before\" string, so use\n    # ReplaceOne for a safe 1:1 swap instead of ReplaceAll.\n    $found = $range.Find.Execute($before, $false, $false, $false, $false, $false, $true, 1, $false, $after, 1)\n\n    if (-not $found) {\n        throw \"Could not find text: $before\"\n    }\n}\n"}
